$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "New--Same Finger": correct the weight for the vf/fv row and
# re-sort the N2:X16 block by the (New Total) column O, ascending.
# -----------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("New--Same Finger")

# Weight correction: vf/fv, mj/jm now costs 3 (was 5)
$ws4.Range("O7").Value = 3

$sortObj = $ws4.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws4.Range("O2:O16"), 0, 1, 0, 0)
$sortObj.SetRange($ws4.Range("N2:X16"))
$sortObj.Header = 2
$sortObj.Apply()

# The sort leaves a stray, now-empty styled cell behind at X7 - clear it
$ws4.Range("X7").Clear()

# Make "New--Same Finger" the active sheet / selected cell
$ws4.Activate()
$ws4.Range("Q19").Select()

# -----------------------------------------------------------------
# Sheet "New--Stretch": drop the now-unused column H, move the
# "cw wc ,o o," row down below the other "Index, Pinky" rows, and
# drop the stale AutoFilter sort-state.
# -----------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("New--Stretch")

# Delete the empty spacer column H (old column I slides left into H, etc.)
$ws5.Columns("H:H").Delete()

$cols = @("A","B","C","D","E","F","G")

# Capture rows 48-52 (values + style) before rewriting them
$rows = @{}
for ($r = 48; $r -le 52; $r++) {
    $rowVals = @{}
    $rowStyles = @{}
    foreach ($col in $cols) {
        $cell = $ws5.Range($col + $r)
        $rowVals[$col] = $cell.Value2
        $rowStyles[$col] = $cell.Style.Name
    }
    $rows[$r] = @{ "vals" = $rowVals; "styles" = $rowStyles }
}

# Move the "cw wc ,o o," entry (originally row 48) down to row 52,
# shifting the four rows below it up by one.
$newOrder = @{48=49; 49=50; 50=51; 51=52; 52=48}
foreach ($destRow in $newOrder.Keys) {
    $srcRow = $newOrder[$destRow]
    foreach ($col in $cols) {
        $destCell = $ws5.Range($col + $destRow)
        $destCell.Value = $rows[$srcRow]["vals"][$col]
        $destCell.Style = $rows[$srcRow]["styles"][$col]
    }
}

# Drop the stale sort-state recorded on the AutoFilter
$ws5.AutoFilterMode = $false
$ws5.Range("A1:G60").AutoFilter()

# Restore this sheet's own selection (it is no longer the active tab)
$ws5.Activate()
$ws5.Range("H10").Select()

# Re-activate "New--Same Finger" last, since it ends up the active tab
$ws4.Activate()
